# Applies the "Optuna Attempt (go back with original)" edit:
# updates forecast numbers on the "Forecast Comparison" sheet and the
# derived totals on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet -------------------------------------------------
# Row 2 (W8)
$wsForecast.Range("D2").Value = 17
$wsForecast.Range("H2").Value = 6.53
$wsForecast.Range("L2").Value = 1.13

# Row 3 (W9)
$wsForecast.Range("D3").Value = 17
$wsForecast.Range("H3").Value = 5.53
$wsForecast.Range("L3").Value = 1.04

# Row 4 (W10)
$wsForecast.Range("H4").Value = 5.92
$wsForecast.Range("L4").Value = 1

# Row 5 (W11)
$wsForecast.Range("H5").Value = 4.1
$wsForecast.Range("L5").Value = 1.09

# Row 6 (W12)
$wsForecast.Range("H6").Value = 3.38
$wsForecast.Range("L6").Value = 0.9399999999999999

# Row 7 (W13)
$wsForecast.Range("H7").Value = 2.62
$wsForecast.Range("L7").Value = 1.11

# Row 8 (W14)
$wsForecast.Range("H8").Value = 1.51
$wsForecast.Range("L8").Value = 1.02

# Row 9 (W15)
$wsForecast.Range("D9").Value = 12
$wsForecast.Range("H9").Value = 0.59
$wsForecast.Range("I9").Value = "Low"
$wsForecast.Range("L9").Value = 1.1

# Row 10 (W16)
$wsForecast.Range("D10").Value = 12
$wsForecast.Range("L10").Value = 0.88

# Row 11 (W17)
$wsForecast.Range("L11").Value = 0.83

# Row 12 (W18)
$wsForecast.Range("L12").Value = 1.19

# Row 13 (W19)
$wsForecast.Range("L13").Value = 0.87

# Row 14 (W20)
$wsForecast.Range("L14").Value = 1.18

# Row 15 (W21)
$wsForecast.Range("L15").Value = 0.8

# Row 16 (W22)
$wsForecast.Range("D16").Value = 13
$wsForecast.Range("L16").Value = 0.96

# Row 17 (W23)
$wsForecast.Range("L17").Value = 1.18

# --- Summary sheet --------------------------------------------------------
# The Value column stores these numbers as plain text, so a leading
# apostrophe forces text entry (matching how Excel keeps numeric-looking
# text as text instead of auto-converting it to a number). Re-applying the
# "Normal" style afterwards drops the quote-prefix formatting flag again so
# the cell format stays exactly as it was before the edit.
function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

Set-TextValue $wsSummary.Range("B9")  "216"
Set-TextValue $wsSummary.Range("B10") "116"
Set-TextValue $wsSummary.Range("B11") "63"
Set-TextValue $wsSummary.Range("B12") "17"
Set-TextValue $wsSummary.Range("B14") "12"
